$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stations")
[void]$ws.Activate()

# The data-cleaning pass removed the blank separator rows that had been left
# between the NSW/QLD, QLD/VIC and VIC/SA blocks. Deleting them (bottom-up,
# so earlier deletions don't shift the row numbers of rows still to be
# deleted) shifts every following block up to close the gaps.
$ws.Rows("98:101").Delete()
$ws.Rows("85:86").Delete()
$ws.Rows("51:54").Delete()

# Excel leaves the selection sitting on the rows that used to be the first
# deleted block (now occupied by the QLD rows that shifted up into them).
[void]$ws.Rows("51:54").Select()
